# "update italian translation #200" — refresh the footer's cached PAGE
# field result (it had gone stale at "1"; the document now paginates to
# page 3 for this footer instance, matching the NUMPAGES field's cached
# "3" total already present: "<PAGE>/<NUMPAGES>").
$d = $word.ActiveDocument

$footer = $d.Sections.First.Footers.Item(1)
$fields = $footer.Range.Fields

for ($i = 1; $i -le $fields.Count; $i++) {
    $field = $fields.Item($i)
    if ($field.Code.Text.Trim() -eq "PAGE") {
        $resultRange = $field.Result
        $resultRange.Find.Execute("1", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "3", 2) | Out-Null
    }
}
